$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F44").Value = 391
$ws.Range("G44").Value = 14236.31
$ws.Range("F67").Value = 183
$ws.Range("G67").Value = 47713.59
$ws.Range("B72").Value = 167741.93
$ws.Range("F102").Value = 14
$ws.Range("G102").Value = 12837.16
$ws.Range("B103").Value = 13994.19
$ws.Range("F112").Value = 179
$ws.Range("G112").Value = 20080.22
$ws.Range("F126").Value = 72
$ws.Range("G126").Value = 9703.440000000001
$ws.Range("B129").Value = 66332.14
$ws.Range("B167").Value = 57756
$ws.Range("E167").Value = 79.37
$ws.Range("F167").Value = -100
$ws.Range("G167").Value = -6644
$ws.Range("B168").Value = 64350
$ws.Range("E168").Value = 70.63
$ws.Range("F168").Value = 2
$ws.Range("G168").Value = 132.88
$ws.Range("F178").Value = 77
$ws.Range("G178").Value = 4891.04
$ws.Range("F180").Value = 35
$ws.Range("G180").Value = 5612.6
$ws.Range("B199").Value = 53851.53
$ws.Range("F202").Value = 0
$ws.Range("G202").Value = 0
$ws.Range("B203").Value = 600.16
$ws.Range("F219").Value = 187
$ws.Range("G219").Value = 23689.16
$ws.Range("F221").Value = 140
$ws.Range("G221").Value = 15720.6
$ws.Range("F223").Value = 3
$ws.Range("G223").Value = 222.9
$ws.Range("B224").Value = 62493.7
$ws.Range("F228").Value = 173
$ws.Range("G228").Value = 3200.5
$ws.Range("F229").Value = 14
$ws.Range("G229").Value = 300.02
$ws.Range("F233").Value = 18
$ws.Range("G233").Value = 2062.8
$ws.Range("B235").Value = 9716.84
$ws.Range("F266").Value = 72
$ws.Range("G266").Value = 6109.92
$ws.Range("F273").Value = 99
$ws.Range("G273").Value = 4205.52
$ws.Range("F283").Value = 0
$ws.Range("G283").Value = 0
$ws.Range("F284").Value = 22
$ws.Range("G284").Value = 2981.88
$ws.Range("F290").Value = 132
$ws.Range("G290").Value = 6186.84
$ws.Range("F291").Value = 0
$ws.Range("G291").Value = 0
$ws.Range("B301").Value = 93829.95
$ws.Range("F308").Value = 34
$ws.Range("G308").Value = 3883.82
$ws.Range("F312").Value = 22
$ws.Range("G312").Value = 3156.56
$ws.Range("B334").Value = -23205.94
$ws.Range("F367").Value = 182
$ws.Range("G367").Value = 25587.38
$ws.Range("B369").Value = 56288.1
$ws.Range("F377").Value = 60
$ws.Range("G377").Value = 9016.200000000001
$ws.Range("B378").Value = 45946.56
$ws.Range("F383").Value = 43
$ws.Range("G383").Value = 41701.83
$ws.Range("B384").Value = 41701.83
$ws.Range("B387").Value = 58047
$ws.Range("D387").Value = 105.54
$ws.Range("E387").Value = 126.1
$ws.Range("F387").Value = 32
$ws.Range("G387").Value = 3377.28
$ws.Range("B388").Value = 47097
$ws.Range("D388").Value = 112.28
$ws.Range("E388").Value = 134.16
$ws.Range("F388").Value = 15
$ws.Range("G388").Value = 1684.2
$ws.Range("F393").Value = 343
$ws.Range("G393").Value = 33133.8
$ws.Range("B395").Value = 48929.29
$ws.Range("F408").Value = 6
$ws.Range("G408").Value = 205.86
$ws.Range("B423").Value = 154344.65
$ws.Range("F438").Value = 48
$ws.Range("G438").Value = 2323.68
$ws.Range("B444").Value = 19554.62
$ws.Range("F461").Value = 30
$ws.Range("G461").Value = 6666.9
$ws.Range("B464").Value = 79371.88
$ws.Range("F525").Value = 341
$ws.Range("G525").Value = 18714.08
$ws.Range("F529").Value = 130
$ws.Range("G529").Value = 11129.3
$ws.Range("B531").Value = 105879.97
$ws.Range("F533").Value = 13
$ws.Range("G533").Value = 430.43
$ws.Range("F535").Value = 99
$ws.Range("G535").Value = 3277.89
$ws.Range("F536").Value = 7
$ws.Range("G536").Value = 302.26
$ws.Range("F540").Value = 102
$ws.Range("G540").Value = 4463.52
$ws.Range("B541").Value = 17982.25
$ws.Range("F552").Value = 34
$ws.Range("G552").Value = 5157.8
$ws.Range("F557").Value = 3
$ws.Range("G557").Value = 2235.57
$ws.Range("B562").Value = 33506.77
$ws.Range("F564").Value = 125
$ws.Range("G564").Value = 15231.25
$ws.Range("B567").Value = 17145.13
$ws.Range("F569").Value = 8
$ws.Range("G569").Value = 1494.88
$ws.Range("F570").Value = 3
$ws.Range("G570").Value = 577.6799999999999
$ws.Range("B579").Value = 11691.38
$ws.Range("F611").Value = 151
$ws.Range("G611").Value = 20098.1
$ws.Range("B613").Value = 20098.1
$ws.Range("F665").Value = 25
$ws.Range("G665").Value = 1338.5
$ws.Range("B674").Value = 9399.120000000001
$ws.Range("F680").Value = 295
$ws.Range("G680").Value = 48117.45
$ws.Range("B686").Value = 49130
$ws.Range("F704").Value = 13
$ws.Range("G704").Value = 2203.89
$ws.Range("B719").Value = 53974.83
$ws.Range("B724").Value = 2121495.95
$ws.Range("B725").Value = 2121495.95
